$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.199962333333334
$ws.Range("H2").Value = 12.599887
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.923694
$ws.Range("N2").Value = 71.77108199999999
$ws.Range("O2").Value = 0.3350100887068322
$ws.Range("P2").Value = 0.3350100887068321
$ws.Range("Q2").Value = 100.4786136741927
$ws.Range("R2").Value = 904.307523067734
$ws.Range("S2").Value = 0.3350100887068322
$ws.Range("T2").Value = 0.3350100887068321

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.199962333333334
$ws.Range("H3").Value = 12.599887
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.063318333333333
$ws.Range("N3").Value = 3.189955
$ws.Range("O3").Value = 0.01488994003909267
$ws.Range("P3").Value = 0.01488994003909267
$ws.Range("Q3").Value = 4.465896948342778
$ws.Range("R3").Value = 40.19307253508501
$ws.Range("S3").Value = 0.01488994003909267
$ws.Range("T3").Value = 0.01488994003909267

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.199962333333334
$ws.Range("H4").Value = 12.599887
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.150636
$ws.Range("N4").Value = 15.451908
$ws.Range("O4").Value = 0.07212577720048599
$ws.Range("P4").Value = 0.07212577720048598
$ws.Range("Q4").Value = 21.63247719271067
$ws.Range("R4").Value = 194.692294734396
$ws.Range("S4").Value = 0.07212577720048599
$ws.Range("T4").Value = 0.07212577720048598

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.199962333333334
$ws.Range("H5").Value = 12.599887
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.265554666666667
$ws.Range("N5").Value = 3.796664
$ws.Range("O5").Value = 0.01772191122087356
$ws.Range("P5").Value = 0.01772191122087356
$ws.Range("Q5").Value = 5.315281930774222
$ws.Range("R5").Value = 47.837537376968
$ws.Range("S5").Value = 0.01772191122087356
$ws.Range("T5").Value = 0.01772191122087356

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.199962333333334
$ws.Range("H6").Value = 12.599887
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.74491
$ws.Range("N6").Value = 17.23473
$ws.Range("O6").Value = 0.08044756000945204
$ws.Range("P6").Value = 0.08044756000945201
$ws.Range("Q6").Value = 24.12840560839
$ws.Range("R6").Value = 217.15565047551
$ws.Range("S6").Value = 0.08044756000945204
$ws.Range("T6").Value = 0.08044756000945201

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.199962333333334
$ws.Range("H7").Value = 12.599887
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 34.26374833333333
$ws.Range("N7").Value = 102.791245
$ws.Range("O7").Value = 0.4798047228232636
$ws.Range("P7").Value = 0.4798047228232635
$ws.Range("Q7").Value = 143.9064523988128
$ws.Range("R7").Value = 1295.158071589315
$ws.Range("S7").Value = 0.4798047228232636
$ws.Range("T7").Value = 0.4798047228232635
